# Weekly data refresh: a new week's price quote was added to the
# "Arveja Verde" (Perfection / $/malla 25 kilos / Provincia de Huasco)
# series. Insert a new row right after the existing row 84 (i.e. at
# row 85), which pushes the former rows 85..158 down to 86..159, then
# populate the newly inserted row with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 85 — shifts rows 85:158 -> 86:159
$ws.Rows.Item(85).Insert()

# Fill in the new row 85 with the latest weekly observation
$ws.Cells.Item(85, 1).Value  = 4
$ws.Cells.Item(85, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(85, 3).Value  = "Los Lagos"
$ws.Cells.Item(85, 4).Value  = 45090
$ws.Cells.Item(85, 5).Value  = 10
$ws.Cells.Item(85, 6).Value  = 100112022
$ws.Cells.Item(85, 7).Value  = "Arveja Verde"
$ws.Cells.Item(85, 8).Value  = "Perfection"
$ws.Cells.Item(85, 9).Value  = "Primera"
$ws.Cells.Item(85, 10).Value = 40
$ws.Cells.Item(85, 11).Value = 42000
$ws.Cells.Item(85, 12).Value = 42000
$ws.Cells.Item(85, 13).Value = 42000
$ws.Cells.Item(85, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(85, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(85, 16).Value = 1680
$ws.Cells.Item(85, 17).Value = 25
$ws.Cells.Item(85, 18).Value = "Hortaliza"

# Preserve the date-formatted style used by the rest of column D
$ws.Cells.Item(85, 4).NumberFormat = $ws.Cells.Item(86, 4).NumberFormat
